# Append three new rows (227-229) to Sheet1, matching the style used by the
# existing data rows (date in column A formatted like the rows above it,
# plain numbers in columns B, C, D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44301, 7, 25, 196.3093835885355),
    @(44302, 2, 27, 212.0141342756183),
    @(44303, 5, 22, 172.7522575579113)
)

$startRow = 227

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Column A: date value, copy style/format from the cell above it
    # (single-cell copy keeps the paste scoped to that one cell).
    $ws.Cells.Item($r - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
